# Add a new "GuilID" property row (guild data module) to the Property sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$row = 11

# Populate in the same left-to-right entry order used when the row was
# authored (this also governs shared-string table append order).
$ws.Cells.Item($row, 10).Value = "工会ID"
$ws.Cells.Item($row, 1).Value = "GuilID"
$ws.Cells.Item($row, 2).Value = "object"
$ws.Cells.Item($row, 3).Value = $true
$ws.Cells.Item($row, 4).Value = $true
$ws.Cells.Item($row, 5).Value = $true
$ws.Cells.Item($row, 6).Value = $true
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = "Friend"

# Match the text-format styling used by the other data rows (A/B/I/J use
# style index 1 -- numFmtId 49, the "@" text format).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 10).NumberFormat = "@"

# The author's selection ended up on E19 after adding the row.
$ws.Range("E19").Select()
